# Update the submission-deadline date on the "Important dates" textbox
# from "February 9, 2018" to "February 14, 2018" (new deadline Feb 14th).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "Important dates" textbox that holds the deadline text.
# Prefer the well-known shape name, but fall back to a text-content
# search so the script still works if shapes get renumbered/renamed.
$targetShape = $null
try {
    $candidate = $s.Shapes.Item("TextBox 17")
    if ($candidate.TextFrame.TextRange.Text -like "*February 9, 2018*") {
        $targetShape = $candidate
    }
} catch {
    $targetShape = $null
}

if ($targetShape -eq $null) {
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -like "*February 9, 2018*") {
                $targetShape = $shp
                break
            }
        }
    }
}

if ($targetShape -ne $null) {
    $tr = $targetShape.TextFrame.TextRange
    $fullText = $tr.Text

    # Find the "9" that immediately follows "February " so only that digit
    # is touched, leaving every other run/character (and its formatting)
    # completely untouched.
    $marker = "February "
    $markerIdx = $fullText.IndexOf($marker)
    if ($markerIdx -ge 0) {
        $digitPos0 = $markerIdx + $marker.Length
        $digitChar = $tr.Characters($digitPos0 + 1, 1)
        if ($digitChar.Text -eq "9") {
            $digitChar.Text = "14"
        }
    }
}
